$d = $word.ActiveDocument

$d.Content.Find.Execute("387×9=", $true, $false, $false, $false, $false, $true, 1, $false, "927×6=", 2) | Out-Null
$d.Content.Find.Execute("259×3=", $true, $false, $false, $false, $false, $true, 1, $false, "558×4=", 2) | Out-Null
$d.Content.Find.Execute("682×5=", $true, $false, $false, $false, $false, $true, 1, $false, "645×7=", 2) | Out-Null
$d.Content.Find.Execute("295×5=", $true, $false, $false, $false, $false, $true, 1, $false, "266×5=", 2) | Out-Null
$d.Content.Find.Execute("176×3=", $true, $false, $false, $false, $false, $true, 1, $false, "553×9=", 2) | Out-Null
$d.Content.Find.Execute("681×9=", $true, $false, $false, $false, $false, $true, 1, $false, "399×2=", 2) | Out-Null
$d.Content.Find.Execute("196×8=", $true, $false, $false, $false, $false, $true, 1, $false, "468×9=", 2) | Out-Null
$d.Content.Find.Execute("972×6=", $true, $false, $false, $false, $false, $true, 1, $false, "739×3=", 2) | Out-Null
$d.Content.Find.Execute("116×3=", $true, $false, $false, $false, $false, $true, 1, $false, "971×9=", 2) | Out-Null
$d.Content.Find.Execute("399×7=", $true, $false, $false, $false, $false, $true, 1, $false, "662×7=", 2) | Out-Null
$d.Content.Find.Execute("313×4=", $true, $false, $false, $false, $false, $true, 1, $false, "563×4=", 2) | Out-Null
$d.Content.Find.Execute("979×2=", $true, $false, $false, $false, $false, $true, 1, $false, "101×6=", 2) | Out-Null
$d.Content.Find.Execute("417×9=", $true, $false, $false, $false, $false, $true, 1, $false, "746×4=", 2) | Out-Null
$d.Content.Find.Execute("401×7=", $true, $false, $false, $false, $false, $true, 1, $false, "843×9=", 2) | Out-Null
$d.Content.Find.Execute("281×6=", $true, $false, $false, $false, $false, $true, 1, $false, "653×6=", 2) | Out-Null
$d.Content.Find.Execute("982×5=", $true, $false, $false, $false, $false, $true, 1, $false, "435×7=", 2) | Out-Null
$d.Content.Find.Execute("853×3=", $true, $false, $false, $false, $false, $true, 1, $false, "872×4=", 2) | Out-Null
$d.Content.Find.Execute("227×4=", $true, $false, $false, $false, $false, $true, 1, $false, "188×4=", 2) | Out-Null
$d.Content.Find.Execute("948×2=", $true, $false, $false, $false, $false, $true, 1, $false, "131×6=", 2) | Out-Null
$d.Content.Find.Execute("599×6=", $true, $false, $false, $false, $false, $true, 1, $false, "434×3=", 2) | Out-Null
$d.Content.Find.Execute("837×8=", $true, $false, $false, $false, $false, $true, 1, $false, "425×4=", 2) | Out-Null
$d.Content.Find.Execute("366×3=", $true, $false, $false, $false, $false, $true, 1, $false, "322×5=", 2) | Out-Null
$d.Content.Find.Execute("935×6=", $true, $false, $false, $false, $false, $true, 1, $false, "827×2=", 2) | Out-Null
$d.Content.Find.Execute("855×2=", $true, $false, $false, $false, $false, $true, 1, $false, "769×2=", 2) | Out-Null
$d.Content.Find.Execute("500×5=", $true, $false, $false, $false, $false, $true, 1, $false, "137×6=", 2) | Out-Null
